$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bayern / Herbstferien (C3): replace footnote text with a date value (31.10.)
$ws.Range("C3").NumberFormat = "d-mmm"
$ws.Range("C3").Value = Get-Date -Year 2024 -Month 10 -Day 31 -Hour 0 -Minute 0 -Second 0

# 2. Hamburg / Himmelf.-Pfingsten (G6): replace footnote reference with plain date text
$ws.Range("G6").Value = "29.05."

# 3. Clear the footnote rows (14-16) and remove their merges
$ws.Range("A14:G14").UnMerge()
$ws.Range("A15:G15").UnMerge()
$ws.Range("A16:G16").UnMerge()

$ws.Range("A14:G14").ClearContents()
$ws.Range("A15:G15").ClearContents()
$ws.Range("A16:G16").ClearContents()

# 4. Update the active selection to match the saved workbook state
$ws.Range("M14").Select()

$wb.Save()
